$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = [char]0x254C
$dashes = "$c$c$c$c$c$c$c"

$ws.Range("A5").Value = $dashes
$ws.Range("B5").Value = $dashes

$ws.Range("B6").Value = 0.06800341606140137
$ws.Range("B7").Value = 0.2119014263153076
$ws.Range("B8").Value = 0.145902156829834
